$wb = $excel.ActiveWorkbook

# Sheet "展览" (Worksheets(1) / sheet1.xml)
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F2").Value = 6689
$ws1.Range("F4").Value = 421
$ws1.Range("F15").Value = 1294
$ws1.Range("F16").Value = 17
$ws1.Range("F17").Value = 3337
$ws1.Range("F18").Value = 18
$ws1.Range("F19").Value = 219
$ws1.Range("F20").Value = 1
$ws1.Range("F21").Value = 1984
$ws1.Range("F22").Value = 102
$ws1.Range("F25").Value = 129

# Sheet "全部类型" (Worksheets(4) / sheet4.xml)
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F2").Value = 6689
$ws4.Range("F4").Value = 421
$ws4.Range("F16").Value = 1294
$ws4.Range("F17").Value = 17
$ws4.Range("F18").Value = 3337
$ws4.Range("F19").Value = 18
$ws4.Range("F20").Value = 219
$ws4.Range("F21").Value = 1
$ws4.Range("F22").Value = 1984
$ws4.Range("F23").Value = 102
$ws4.Range("F26").Value = 129
